$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Latitude/Longitude for row 2 (Blackstone Hall)
$ws.Range("G2").Value = 42.252001
$ws.Range("H2").Value = -71.82108

# Add Latitude/Longitude for row 3
$ws.Range("G3").Value = 42.25001
$ws.Range("H3").Value = -71.826244

# Update the sheet view: remove the scrolled topLeftCell and update selection
$ws.Range("H4").Select()
